# Rename the second worksheet ("x1" -> "Word") to match the new vocabulary topic
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Name = "Word"

# Overwrite the English/Thai word-pair table with the new vocabulary list
# (24 rows x 2 columns -- replaces the previous business-vocabulary content)
$ws.Cells.Item(1, 1).Value = 'Hello'
$ws.Cells.Item(1, 2).Value = 'สวัสดี'
$ws.Cells.Item(2, 1).Value = 'Thank you'
$ws.Cells.Item(2, 2).Value = 'ขอบคุณ'
$ws.Cells.Item(3, 1).Value = 'Rice'
$ws.Cells.Item(3, 2).Value = 'ข้าว'
$ws.Cells.Item(4, 1).Value = 'Water'
$ws.Cells.Item(4, 2).Value = 'น้ำ'
$ws.Cells.Item(5, 1).Value = 'Book'
$ws.Cells.Item(5, 2).Value = 'หนังสือ'
$ws.Cells.Item(6, 1).Value = 'Television'
$ws.Cells.Item(6, 2).Value = 'โทรทัศน์'
$ws.Cells.Item(7, 1).Value = 'Jewelry'
$ws.Cells.Item(7, 2).Value = 'เครื่องประดับ'
$ws.Cells.Item(8, 1).Value = 'Animal'
$ws.Cells.Item(8, 2).Value = 'สัตว์'
$ws.Cells.Item(9, 1).Value = 'Car'
$ws.Cells.Item(9, 2).Value = 'รถยนต์'
$ws.Cells.Item(10, 1).Value = 'Pen'
$ws.Cells.Item(10, 2).Value = 'ปากกา'
$ws.Cells.Item(11, 1).Value = 'House'
$ws.Cells.Item(11, 2).Value = 'บ้าน'
$ws.Cells.Item(12, 1).Value = 'Tree'
$ws.Cells.Item(12, 2).Value = 'ต้นไม้'
$ws.Cells.Item(13, 1).Value = 'Flower'
$ws.Cells.Item(13, 2).Value = 'ดอกไม้'
$ws.Cells.Item(14, 1).Value = 'Food'
$ws.Cells.Item(14, 2).Value = 'อาหาร'
$ws.Cells.Item(15, 1).Value = 'Fruit'
$ws.Cells.Item(15, 2).Value = 'ผลไม้'
$ws.Cells.Item(16, 1).Value = 'Vegetable'
$ws.Cells.Item(16, 2).Value = 'ผัก'
$ws.Cells.Item(17, 1).Value = 'Mother'
$ws.Cells.Item(17, 2).Value = 'แม่'
$ws.Cells.Item(18, 1).Value = 'Father'
$ws.Cells.Item(18, 2).Value = 'พ่อ'
$ws.Cells.Item(19, 1).Value = 'Brother'
$ws.Cells.Item(19, 2).Value = 'พี่ชาย'
$ws.Cells.Item(20, 1).Value = 'Sister'
$ws.Cells.Item(20, 2).Value = 'พี่สาว'
$ws.Cells.Item(21, 1).Value = 'Friend'
$ws.Cells.Item(21, 2).Value = 'เพื่อน'
$ws.Cells.Item(22, 1).Value = 'Love'
$ws.Cells.Item(22, 2).Value = 'รัก'
$ws.Cells.Item(23, 1).Value = 'Money'
$ws.Cells.Item(23, 2).Value = 'เงิน'
$ws.Cells.Item(24, 1).Value = 'Time'
$ws.Cells.Item(24, 2).Value = 'เวลา'

# Move the active selection to D19 (matches the saved view state)
[void]$ws.Range("D19").Select()

# Give the sheet a defined print/page setup (adds <pageSetup orientation="portrait"/>)
$ws.PageSetup.Orientation = 1

